$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 13, shifting existing rows 13-14 down to 14-15.
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the latest week's data.
$ws.Cells.Item(13, 1).Value = 1
$ws.Cells.Item(13, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(13, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(13, 4).Value = 45205
$ws.Cells.Item(13, 5).Value = 15
$ws.Cells.Item(13, 6).Value = 100112030
$ws.Cells.Item(13, 7).Value = "Poroto granado"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 3500
$ws.Cells.Item(13, 11).Value = 1400
$ws.Cells.Item(13, 12).Value = 1500
$ws.Cells.Item(13, 13).Value = 1457
$ws.Cells.Item(13, 14).Value = "$/kilo"
$ws.Cells.Item(13, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(13, 16).Value = 1457
$ws.Cells.Item(13, 17).Value = 1
$ws.Cells.Item(13, 18).Value = "Hortaliza"
